$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item('LP1912')
$ws2 = $wb.Worksheets.Item('LP1912-215')
$ws3 = $wb.Worksheets.Item('6203-6173')

# --- Sheet 1: LP1912 ---
$ws1.Range('A2').Value = 'Última actualización: 11:49:24'
$ws1.Range('A3').Value = 'Total filas: 203'
$ws1.Range('A88').Value = '07:28:14'
$ws1.Range('C88').Value = '16_P MOR-SANTA ANA'
$ws1.Range('D88').Value = 55
$ws1.Range('A89').Value = '08:13:38'
$ws1.Range('C89').Value = '215B_EL PATO'
$ws1.Range('D89').Value = 10
$ws1.Range('A103').Value = '08:56:26'
$ws1.Range('C103').Value = '23_HERNANDEZ'
$ws1.Range('D103').Value = 8
$ws1.Range('A104').Value = '08:48:01'
$ws1.Range('C104').Value = '11_ETCHEVERRY'
$ws1.Range('D104').Value = 16
$ws1.Range('C163').Value = '26_HERNANDEZ'
$ws1.Range('C164').Value = '16_SANTA ANA'
$ws1.Range('A175').Value = '11:49:23'
$ws1.Range('D175').Value = 2
$ws1.Range('A177').Value = '11:49:23'
$ws1.Range('D177').Value = 10
$ws1.Range('A178').Value = '11:49:23'
$ws1.Range('D178').Value = 13
$ws1.Range('A179').Value = '11:49:23'
$ws1.Range('D179').Value = 15
$ws1.Range('A181').Value = '11:49:23'
$ws1.Range('D181').Value = 17
$ws1.Range('A182').Value = '11:49:23'
$ws1.Range('D182').Value = 24
$ws1.Range('A183').Value = '11:49:23'
$ws1.Range('C183').Value = '17_ROMERO'
$ws1.Range('D183').Value = 25
$ws1.Range('A184').Value = '10:57:35'
$ws1.Range('C184').Value = '10_OLMOS'
$ws1.Range('D184').Value = 77
$ws1.Range('A185').Value = '11:49:23'
$ws1.Range('B185').Value = '12:16'
$ws1.Range('C185').Value = '16_SANTA ANA'
$ws1.Range('D185').Value = 27
$ws1.Range('A186').Value = '10:27:29'
$ws1.Range('B186').Value = '12:17'
$ws1.Range('D186').Value = 110
$ws1.Range('A187').Value = '11:49:23'
$ws1.Range('C187').Value = '14_ABASTO'
$ws1.Range('D187').Value = 31
$ws1.Range('A188').Value = '11:49:23'
$ws1.Range('B188').Value = '12:20'
$ws1.Range('C188').Value = '215A_EL PATO'
$ws1.Range('D188').Value = 31
$ws1.Range('A189').Value = '11:49:23'
$ws1.Range('B189').Value = '12:21'
$ws1.Range('C189').Value = '26_HERNANDEZ'
$ws1.Range('D189').Value = 32
$ws1.Range('A190').Value = '11:49:23'
$ws1.Range('B190').Value = '12:34'
$ws1.Range('C190').Value = '11_ETCHEVERRY'
$ws1.Range('D190').Value = 45
$ws1.Range('A191').Value = '11:49:23'
$ws1.Range('B191').Value = '12:34'
$ws1.Range('C191').Value = '23_HERNANDEZ'
$ws1.Range('D191').Value = 45
$ws1.Range('B192').Value = '12:36'
$ws1.Range('C192').Value = '27_EL RETIRO'
$ws1.Range('D192').Value = 74
$ws1.Range('A193').Value = '11:49:23'
$ws1.Range('B193').Value = '12:37'
$ws1.Range('C193').Value = '27_EL RETIRO'
$ws1.Range('D193').Value = 48
$ws1.Range('A194').Value = '11:49:23'
$ws1.Range('B194').Value = '12:38'
$ws1.Range('C194').Value = '17_179 Y 38'
$ws1.Range('D194').Value = 49
$ws1.Range('A195').Value = '11:49:23'
$ws1.Range('B195').Value = '12:41'
$ws1.Range('C195').Value = '10_OLMOS'
$ws1.Range('D195').Value = 52
$ws1.Range('A196').Value = '11:49:23'
$ws1.Range('B196').Value = '12:45'
$ws1.Range('C196').Value = '16_SANTA ANA'
$ws1.Range('D196').Value = 56
$ws1.Range('A197').Value = '11:49:23'
$ws1.Range('B197').Value = '12:48'
$ws1.Range('C197').Value = '11_ETCHEVERRY'
$ws1.Range('D197').Value = 59
$ws1.Range('E197').Value = 'LP1912'
$ws1.Range('A198').Value = '11:22:44'
$ws1.Range('B198').Value = '12:50'
$ws1.Range('C198').Value = '15_ABASTO'
$ws1.Range('D198').Value = 88
$ws1.Range('E198').Value = 'LP1912'
$ws1.Range('A199').Value = '11:49:23'
$ws1.Range('B199').Value = '13:02'
$ws1.Range('C199').Value = '15_ABASTO'
$ws1.Range('D199').Value = 73
$ws1.Range('E199').Value = 'LP1912'
$ws1.Range('A200').Value = '11:49:23'
$ws1.Range('B200').Value = '13:06'
$ws1.Range('C200').Value = '16_P MOR-SANTA ANA'
$ws1.Range('D200').Value = 77
$ws1.Range('E200').Value = 'LP1912'
$ws1.Range('A201').Value = '11:22:44'
$ws1.Range('B201').Value = '13:10'
$ws1.Range('C201').Value = '10_OLMOS'
$ws1.Range('D201').Value = 108
$ws1.Range('E201').Value = 'LP1912'
$ws1.Range('A202').Value = '11:22:44'
$ws1.Range('B202').Value = '13:13'
$ws1.Range('C202').Value = '215D_EL PATO'
$ws1.Range('D202').Value = 111
$ws1.Range('E202').Value = 'LP1912'
$ws1.Range('A203').Value = '11:49:23'
$ws1.Range('B203').Value = '13:14'
$ws1.Range('C203').Value = '215D_EL PATO'
$ws1.Range('D203').Value = 85
$ws1.Range('E203').Value = 'LP1912'
$ws1.Range('A204').Value = '11:49:23'
$ws1.Range('B204').Value = '13:20'
$ws1.Range('C204').Value = '10_OLMOS'
$ws1.Range('D204').Value = 91
$ws1.Range('E204').Value = 'LP1912'
$ws1.Range('A205').Value = '11:49:23'
$ws1.Range('B205').Value = '13:21'
$ws1.Range('C205').Value = '26_HERNANDEZ'
$ws1.Range('D205').Value = 92
$ws1.Range('E205').Value = 'LP1912'
$ws1.Range('A206').Value = '11:49:23'
$ws1.Range('B206').Value = '13:26'
$ws1.Range('C206').Value = '15_ABASTO'
$ws1.Range('D206').Value = 97
$ws1.Range('E206').Value = 'LP1912'
$ws1.Range('A207').Value = '11:49:23'
$ws1.Range('B207').Value = '13:26'
$ws1.Range('C207').Value = '14_ABASTO'
$ws1.Range('D207').Value = 97
$ws1.Range('E207').Value = 'LP1912'
$ws1.Range('A208').Value = '11:49:23'
$ws1.Range('B208').Value = '13:46'
$ws1.Range('C208').Value = '17_ROMERO'
$ws1.Range('D208').Value = 117
$ws1.Range('E208').Value = 'LP1912'

# --- Sheet 2: LP1912-215 ---
$ws2.Range('A2').Value = 'Última actualización: 11:49:24'
$ws2.Range('A3').Value = 'Total filas: 27'
$ws2.Range('A29').Value = '11:49:23'
$ws2.Range('D29').Value = 2
$ws2.Range('A30').Value = '11:49:23'
$ws2.Range('D30').Value = 31
$ws2.Range('A32').Value = '11:49:23'
$ws2.Range('B32').Value = '13:14'
$ws2.Range('C32').Value = '215D_EL PATO'
$ws2.Range('D32').Value = 85
$ws2.Range('E32').Value = 'LP1912'

# --- Sheet 3: 6203-6173 ---
$ws3.Range('A2').Value = 'Última actualización: 11:49:24'
$ws3.Range('A3').Value = 'Total filas: 32'
$ws3.Range('A34').Value = '11:49:23'
$ws3.Range('D34').Value = 15
$ws3.Range('A36').Value = '11:49:23'
$ws3.Range('B36').Value = '12:54'
$ws3.Range('C36').Value = '215C_LA PLATA'
$ws3.Range('D36').Value = 65
$ws3.Range('E36').Value = 'L6203'
$ws3.Range('A37').Value = '11:49:23'
$ws3.Range('B37').Value = '13:31'
$ws3.Range('C37').Value = '215B_LP-P MOR-1 Y 57'
$ws3.Range('D37').Value = 102
$ws3.Range('E37').Value = 'L6173'
